$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.723.18"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "2.340.61"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "501.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.37%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -3.08%  "
$ws.Range("D9").Value = "2.348.12"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").Value = "2.757.32"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").Value = "55.710.44"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("D18").Value = "2.359.72"
$ws.Range("E18").Value = "  -3.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "309.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").Value = "0.0₃0701"
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -5.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("E38").Value = "  -4.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.816"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("E41").Value = "  -4.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "126.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.77%  "
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.553"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("E46").Value = "  -2.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "236.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.74%  "
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("E49").Value = "  -2.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  +0.08%  "
